$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The StatQuery cell (B2) was missing a trailing space after "LIMIT 100",
# which was causing the downstream reader to truncate/limit results
# instead of reading all tabs. Append the missing trailing space.
$cell = $ws.Range("B2")
$cell.Value = $cell.Text + " "

# Update the active selection on the sheet to C2 (as last saved by the author).
$ws.Range("C2").Select()
